$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 23
$ws.Range("H4").Value = 6.4
$ws.Range("I4").Value = 1.12
$ws.Range("N4").Value = 1.5
$ws.Range("O4").Value = 2.27
$ws.Range("R4").Value = 2.55
$ws.Range("S4").Value = 1.39
$ws.Range("T4").Value = 60
$ws.Range("U4").Value = 300
$ws.Range("V4").Value = 90
$ws.Range("X4").Value = 600
$ws.Range("Y4").Value = 350
$ws.Range("Z4").Value = 13.5
$ws.Range("AA4").Value = 15.5
$ws.Range("AB4").Value = 45
$ws.Range("AD4").Value = 7
$ws.Range("AF4").Value = 11.25
$ws.Range("AG4").Value = 5.7
$ws.Range("AI4").Value = 50
$ws.Range("G6").Value = 1.22
$ws.Range("H6").Value = 5.4
$ws.Range("I6").Value = 11.25
$ws.Range("L6").Value = 1.19
$ws.Range("M6").Value = 3.7
$ws.Range("N6").Value = 1.6
$ws.Range("O6").Value = 2.07
$ws.Range("R6").Value = 2.22
$ws.Range("S6").Value = 1.52
$ws.Range("T6").Value = 6.8
$ws.Range("U6").Value = 5.6
$ws.Range("V6").Value = 9.5
$ws.Range("W6").Value = 6.7
$ws.Range("Y6").Value = 37
$ws.Range("Z6").Value = 13
$ws.Range("AA6").Value = 11.5
$ws.Range("AB6").Value = 30
$ws.Range("AC6").Value = 175
$ws.Range("AD6").Value = 27
$ws.Range("AE6").Value = 90
$ws.Range("AF6").Value = 37
$ws.Range("AG6").Value = 400
$ws.Range("AH6").Value = 175
$ws.Range("AI6").Value = 150
$ws.Range("G7").Value = 1.18
$ws.Range("H7").Value = 5.8
$ws.Range("I7").Value = 17.5
$ws.Range("K7").Value = 8
$ws.Range("L7").Value = 1.24
$ws.Range("M7").Value = 3.65
$ws.Range("N7").Value = 1.72
$ws.Range("O7").Value = 2
$ws.Range("P7").Value = 1.36
$ws.Range("Q7").Value = 2.9
$ws.Range("R7").Value = 2.87
$ws.Range("S7").Value = 1.37
$ws.Range("T7").Value = 5.7
$ws.Range("U7").Value = 4.75
$ws.Range("V7").Value = 11
$ws.Range("W7").Value = 5.8
$ws.Range("Y7").Value = 55
$ws.Range("Z7").Value = 8
$ws.Range("AA7").Value = 13
$ws.Range("AB7").Value = 45
$ws.Range("AC7").Value = 350
$ws.Range("AD7").Value = 35
$ws.Range("AE7").Value = 150
$ws.Range("AF7").Value = 60
$ws.Range("AH7").Value = 400
$ws.Range("AI7").Value = 300
$ws.Range("G10").Value = 2.25
$ws.Range("H10").Value = 2.82
$ws.Range("I10").Value = 3.7
$ws.Range("J10").Value = 1.08
$ws.Range("K10").Value = 6.8
$ws.Range("L10").Value = 1.33
$ws.Range("M10").Value = 3.15
$ws.Range("N10").Value = 1.98
$ws.Range("O10").Value = 1.78
$ws.Range("P10").Value = 1.4
$ws.Range("Q10").Value = 2.82
$ws.Range("R10").Value = 1.65
$ws.Range("S10").Value = 2.1
$ws.Range("T10").Value = 7.5
$ws.Range("U10").Value = 12.5
$ws.Range("X10").Value = 19.5
$ws.Range("Y10").Value = 28
$ws.Range("Z10").Value = 6.8
$ws.Range("AB10").Value = 12.5
$ws.Range("AC10").Value = 55
$ws.Range("AD10").Value = 10.5
$ws.Range("AG10").Value = 70
$ws.Range("AI10").Value = 37
$ws.Range("AJ10").Value = 450
$ws.Range("H14").Value = 3.1
$ws.Range("I14").Value = 3
$ws.Range("J14").Value = 1.1
$ws.Range("K14").Value = 7
$ws.Range("L14").Value = 1.44
$ws.Range("M14").Value = 2.63
$ws.Range("Z14").Value = 7
$ws.Range("N18").Value = 2
$ws.Range("O18").Value = 1.8
$ws.Range("P29").Value = 1.53
$ws.Range("R29").Value = 1.93
$ws.Range("T29").Value = 6.3
$ws.Range("Z29").Value = 6.9
$ws.Range("AD29").Value = 7.7
$ws.Range("AG29").Value = 45
$ws.Range("AI29").Value = 45
$ws.Range("G32").Value = 2.45
$ws.Range("I32").Value = 2.4
$ws.Range("W32").Value = 26
$ws.Range("AD32").Value = 11
$ws.Range("AH32").Value = 19
$ws.Range("G37").Value = 2.35
$ws.Range("L37").Value = 1.39
$ws.Range("M37").Value = 2.55
$ws.Range("N37").Value = 2.12
$ws.Range("P37").Value = 1.45
$ws.Range("Q37").Value = 2.37
$ws.Range("R37").Value = 1.85
$ws.Range("S37").Value = 1.75
$ws.Range("Z37").Value = 7.5
$ws.Range("AA37").Value = 5.8
$ws.Range("AB37").Value = 15.5
$ws.Range("AD37").Value = 7.8
$ws.Range("AE37").Value = 14.5
$ws.Range("AF37").Value = 11.25
$ws.Range("AI37").Value = 45
$ws.Range("AJ37").Value = 800
